# Add a "Save" column (H) to the s_vals sheet, mirroring the header style
# used by the existing columns and filling in the per-row save flag.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting of the last existing header cell (G1) onto the new
# header cell (H1) so it keeps the bold / centered / bordered style, then
# set its text.
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122) # xlPasteFormats
$ws.Range("H1").Value = "Save"

# Fill in the "Save" values for each data row (2-13).
$saveValues = @(1, 0, 1, 1, 1, 1, 1, 1, 1, 1, 1, 1)
for ($i = 0; $i -lt $saveValues.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 8).Value = $saveValues[$i]
}
